# Capitalizes the month labels in column A (rows 2-22, e.g. "jan/23" -> "Jan/23")
# on both the "DEC" and "FEC" sheets, removes a stray empty formatted row left
# over at the bottom of "FEC", and switches the active sheet/selection to FEC.
$wb = $excel.ActiveWorkbook

function Capitalize-Months {
    param($ws)
    for ($r = 2; $r -le 22; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        $val = $cell.Value2
        if ($val -ne $null -and $val.Length -gt 0) {
            $cell.Value = $val.Substring(0,1).ToUpper() + $val.Substring(1)
        }
    }
}

$wsDec = $wb.Worksheets.Item("DEC")
$wsFec = $wb.Worksheets.Item("FEC")

Capitalize-Months $wsDec
Capitalize-Months $wsFec

# Remove leftover empty formatted row 23 on FEC
$wsFec.Rows.Item(23).Delete()

# Selections / active sheet per target diff
$wsDec.Range("A2:A22").Select()
$wsFec.Activate()
$wsFec.Range("D9").Select()
